$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header cell text to reflect the new "through" date
$ws.Name = "Through 2022-07-09"
$ws.Range("B1").Value = "July 2022 (through July 09)"

# Helper to set a cell's value directly
function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Updates to existing values (increment counts for newly added incidents)
Set-Cell "B2"  2   # Austin, July 2022
Set-Cell "B3"  3   # Englewood, July 2022
Set-Cell "P3"  2   # Englewood, July 2020
Set-Cell "AK7" 2   # Roseland, July 2017
Set-Cell "I8"  4   # North Lawndale, July 2021
Set-Cell "AK39" 2  # Wicker Park, July 2017
Set-Cell "P52" 4   # Chatham, July 2020
Set-Cell "I78" 2   # Lake View, July 2021

# New cells that previously had no value
Set-Cell "W4"  1   # Auburn Gresham, July 2019
Set-Cell "W5"  1   # Garfield Park, July 2019
Set-Cell "I6"  1   # Grand Crossing, July 2021
Set-Cell "AR17" 1  # Brighton Park, July 2016
Set-Cell "P19" 1   # South Shore, July 2020
Set-Cell "AK19" 1  # South Shore, July 2017
Set-Cell "P23" 1   # South Chicago, July 2020
Set-Cell "I26" 1   # Little Village, July 2021
Set-Cell "I34" 1   # Riverdale, July 2021
Set-Cell "AD38" 1  # West Town, July 2018
Set-Cell "P56" 1   # Belmont Cragin, July 2020
Set-Cell "AD82" 1  # Morgan Park, July 2018
